# Applies the "added environment & reworked optimisation" schedule rework:
#   - introduces a new "cleanup" activity (new shared string + new theme fill/style)
#   - reshuffles several 15-minute schedule slots between any/eat/relax/cleanup

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) B36: "any" -> "eat"  (reuse the existing red "eat" fill/style, e.g. B34)
# ---------------------------------------------------------------------------
$ws.Range("B36").Value = "eat"
$ws.Range("B34").Copy()
$ws.Range("B36").PasteSpecial(-4122)   # xlPasteFormats

# ---------------------------------------------------------------------------
# 2) Cells that become "relax" (reuse the existing "relax" fill/style, e.g. E39)
# ---------------------------------------------------------------------------
$relaxCells = @(
    "E50","E51","E52","E53",
    "B56","E56","B57","E57","E58","E59","E60","E61","E62",
    "B90","E90","B91","E91","B92","E92"
)
foreach ($ref in $relaxCells) {
    $ws.Range($ref).Value = "relax"
}
$ws.Range("E39").Copy()
foreach ($ref in $relaxCells) {
    $ws.Range($ref).PasteSpecial(-4122)   # xlPasteFormats
}

# ---------------------------------------------------------------------------
# 3) Cells that become "cleanup" (brand-new activity/style - theme accent4 fill)
# ---------------------------------------------------------------------------
$cleanupCells = @(
    "E47","E48","E49",
    "E71","E72","E73","E74","E75","E76","E77",
    "B82","B83"
)
foreach ($ref in $cleanupCells) {
    $ws.Range($ref).Value = "cleanup"
}
# Build the new fill on the first cell, then propagate it to the rest.
$first = $cleanupCells[0]
$ws.Range($first).Interior.ThemeColor = 8
$ws.Range($first).Copy()
foreach ($ref in $cleanupCells) {
    $ws.Range($ref).PasteSpecial(-4122)   # xlPasteFormats
}

$excel.CutCopyMode = 0
